$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the existing "sum" header (G1) onto the new H1
# header cell so the new column matches the look of the rest of the header
# row (bold, bordered, centered) without minting a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column for each data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
